$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (B2) is corrected: the trailing `Cohort` output
# column (and the now-unneeded `cohort` match feeding it) is dropped from
# the Cypher text, so the query's last returned column becomes
# "Response to Treatment".
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
 WHERE labels(parent)[0] IN ["diagnosis"] 
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newQuery
$ws.Rows.Item(2).AutoFit()

# Mirror the author's window/view state: scrolled back to the top of the
# sheet with row 2 / B2 selected, at 100% zoom.
$ws.Activate()
$ws.Range("A2").Select()
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 100

$wb.Save()
